# "excel folder been created" - append 3 new form submissions (rows 21-23)
# to the Submissions sheet.
#
# The sheet stores every value as text (shared strings) - even the
# "numeric-looking" Age/Phone columns - so whenever we write one of those
# columns we briefly force a text number-format on that single cell, assign
# the value, and then clear the format again so no stray per-cell style is
# left behind. Each row is filled in strict column order (A -> G) so the
# shared-string table grows in the same order Excel would naturally create
# it while typing the row left-to-right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# --- Row 21 : Akhil, 36 ----------------------------------------------------
$ws.Range("A21").Value = "8/9/2025, 6:15:43 pm"
$ws.Range("B21").Value = "Akhil"
Set-TextValue "C21" "36"
$ws.Range("D21").Value = "Male"
Set-TextValue "E21" "8008936289"
$ws.Range("F21").Value = "digitalsales@progenicslabs.com"
$ws.Range("G21").Value = '{"name":"Akhil","age":"36","sex":"Male","phone":"8008936289","email":"digitalsales@progenicslabs.com","sec-1":["Bloating or abdominal distension","Change in bowel frequency and stool form and shape"],"sec-2":["Visceral hypersensitivity","Dysregulated gut motility"],"sec-3":["Abdominal pain >1 day/week","Pain related to defecation"],"sec-4":["Fecal Incontinence","Weight loss","Fever","Nocturnal symptoms"]}'

# --- Row 22 : udaykiran, 23 --------------------------------------------------
$ws.Range("A22").Value = "8/9/2025, 6:22:44 pm"
$ws.Range("B22").Value = "udaykiran"
Set-TextValue "C22" "23"
$ws.Range("D22").Value = "Male"
Set-TextValue "E22" "9966960202"
$ws.Range("F22").Value = "gogulaudaykiran2204@gmail.com"
$ws.Range("G22").Value = '{"name":"udaykiran","age":"23","sex":"Male","phone":"9966960202","email":"gogulaudaykiran2204@gmail.com","sec-1":["Bloating or abdominal distension","Change in bowel frequency and stool form and shape"],"sec-2":["Abnormal pain signaling","Microbiota disturbance","Gut Inflammation"],"sec-3":["Abdominal pain >1 day/week","Pain related to defecation","Type 2: Lumpy and sausage-shaped","Type 5: Soft blobs with clear-cut edges","Type 7: Watery, no solid pieces"],"sec-4":["Blood in stool","Anemia or low hemoglobin","Abdominal mass","Fecal Incontinence","Nocturnal symptoms"]}'

# --- Row 23 : uday kiran, 23 -------------------------------------------------
$ws.Range("A23").Value = "8/9/2025, 6:32:13 pm"
$ws.Range("B23").Value = "uday kiran"
Set-TextValue "C23" "23"
$ws.Range("D23").Value = "Male"
Set-TextValue "E23" "9966960202"
$ws.Range("F23").Value = "gogulaudaykiran2204@gmail.com"
$ws.Range("G23").Value = '{"name":"uday kiran","age":"23","sex":"Male","country-code":"+91","phone":"9966960202","email":"gogulaudaykiran2204@gmail.com","sec-1":["Recurrent abdominal pain (>1 day/week in the last 3 months)","Bloating or abdominal distension","Change in bowel frequency and stool form and shape"],"sec-2":["Abnormal pain signaling","Dysregulated gut motility","Microbiota disturbance","Gut Inflammation","Dietary trigger"],"sec-3":["Symptoms > 6 months in duration","Abdominal pain >1 day/week","Pain related to defecation","Type 2: Lumpy and sausage-shaped","Type 4: Smooth, soft sausage or snake","Type 5: Soft blobs with clear-cut edges","Type 6: Mushy with ragged edges","Type 7: Watery, no solid pieces"],"sec-4":["Abdominal mass","Fecal Incontinence","Weight loss","Fever","Nocturnal symptoms","Family history of IBD, celiac, cancer","New symptom onset (<6 months)","Recent antibiotic use","Extra-intestinal signs (rash, arthritis, eye Inflammation)"]}'
